$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.783420205116272
$ws.Range("B1").Value = 3.262901782989502
$ws.Range("C1").Value = 2.574892044067383
$ws.Range("D1").Value = 2.46576452255249
$ws.Range("E1").Value = 2.316381454467773
